$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(24,25,25,24,24,24)
    3 = @(24,25,25,25,25,24)
    4 = @(25,25,24,25,25,24)
    5 = @(26,26,27,27,25,24)
    6 = @(25,25,27,26,25,25)
    7 = @(24,25,26,27,25,24)
    8 = @(26,25,25,26,25,24)
    9 = @(25,23,25,24,23,23)
}

$cols = @("B","C","D","E","F","G")

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $rowVals[$i]
        $ws.Range("$col$row").Value = $val
    }
}
